$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated corona report numbers (#53).
# Column A (Bundesland names) is unchanged; columns B-H are refreshed.

$ws.Range("B2").Value = "0,8"
$ws.Range("C2").Value = 157
$ws.Range("D2").Value = 309
$ws.Range("E2").Value = 688782
$ws.Range("F2").Value = 2.35087810554394
$ws.Range("G2").Value = 32.9869659348083
$ws.Range("H2").Value = 2.92784804444886

$ws.Range("B3").Value = "1,3"
$ws.Range("C3").Value = 94
$ws.Range("D3").Value = 210
$ws.Range("E3").Value = 38426
$ws.Range("F3").Value = 1.00651807148837
$ws.Range("G3").Value = 47.0602007667769
$ws.Range("H3").Value = 5.45319708366641

$ws.Range("B4").Value = "0,8"
$ws.Range("C4").Value = 119
$ws.Range("D4").Value = 236
$ws.Range("E4").Value = 14273
$ws.Range("F4").Value = 2.23770106206351
$ws.Range("G4").Value = 23.3942326055128
$ws.Range("H4").Value = 4.59014567027825

$ws.Range("B5").Value = "0,7"
$ws.Range("C5").Value = 113
$ws.Range("D5").Value = 212
$ws.Range("E5").Value = 54168
$ws.Range("F5").Value = 1.53656521560727
$ws.Range("G5").Value = 31.1758279422475
$ws.Range("H5").Value = 2.12277764950358

$ws.Range("B6").Value = "1,0"
$ws.Range("C6").Value = 85
$ws.Range("D6").Value = 81
$ws.Range("E6").Value = 6999
$ws.Range("F6").Value = 2.12580115736595
$ws.Range("G6").Value = 45.9946146078761
$ws.Range("H6").Value = 2.35710206054216

$ws.Range("B7").Value = "0,7"
$ws.Range("C7").Value = 142
$ws.Range("D7").Value = 248
$ws.Range("E7").Value = 124637
$ws.Range("F7").Value = 2.42648708677516
$ws.Range("G7").Value = 41.2788389295771
$ws.Range("H7").Value = 0

$ws.Range("B8").Value = "0,9"
$ws.Range("C8").Value = 150
$ws.Range("D8").Value = 339
$ws.Range("E8").Value = 57475
$ws.Range("F8").Value = 2.4300263355428
$ws.Range("G8").Value = 40.0047151847083
$ws.Range("H8").Value = 3.79047567190296

$ws.Range("B9").Value = "0,9"
$ws.Range("C9").Value = 133
$ws.Range("D9").Value = 286
$ws.Range("E9").Value = 36466
$ws.Range("F9").Value = 1.99733603849432
$ws.Range("G9").Value = 27.825295221278
$ws.Range("H9").Value = 2.51567198048448

$ws.Range("B10").Value = "0,7"
$ws.Range("C10").Value = 129
$ws.Range("D10").Value = 214
$ws.Range("E10").Value = 73776
$ws.Range("F10").Value = 2.38435680751512
$ws.Range("G10").Value = 12.4717826742867
$ws.Range("H10").Value = 5.32842316333814

$ws.Range("B11").Value = "1,0"
$ws.Range("C11").Value = 149
$ws.Range("D11").Value = 291
$ws.Range("E11").Value = 135986
$ws.Range("F11").Value = 2.75576569648596
$ws.Range("G11").Value = 39.7244209909118
$ws.Range("H11").Value = 3.8101068071038

$ws.Range("B12").Value = "1,0"
$ws.Range("C12").Value = 182
$ws.Range("D12").Value = 388
$ws.Range("E12").Value = 10130
$ws.Range("F12").Value = 2.26358235542671
$ws.Range("G12").Value = 30.9951775537045
$ws.Range("H12").Value = 10.1009301405106

$ws.Range("B13").Value = "0,9"
$ws.Range("C13").Value = 191
$ws.Range("D13").Value = 438
$ws.Range("E13").Value = 33798
$ws.Range("F13").Value = 2.958012432787
$ws.Range("G13").Value = 70.8310210876804
$ws.Range("H13").Value = 10.7154861831111

$ws.Range("B14").Value = "0,6"
$ws.Range("C14").Value = 260
$ws.Range("D14").Value = 456
$ws.Range("E14").Value = 15091
$ws.Range("F14").Value = 2.10353888923915
$ws.Range("G14").Value = 4.78543777283585
$ws.Range("H14").Value = 1.22931298328823

$ws.Range("B15").Value = "1,7"
$ws.Range("C15").Value = 122
$ws.Range("D15").Value = 243
$ws.Range("E15").Value = 26926
$ws.Range("F15").Value = 0.952343642150114
$ws.Range("G15").Value = 64.4183977813825
$ws.Range("H15").Value = 1.2598276025386

$ws.Range("B16").Value = "0,6"
$ws.Range("C16").Value = 309
$ws.Range("D16").Value = 633
$ws.Range("E16").Value = 23701
$ws.Range("F16").Value = 3.8995366126134
$ws.Range("G16").Value = 9.00053681134064
$ws.Range("H16").Value = 0.517552335867094

$ws.Range("B17").Value = "1,1"
$ws.Range("C17").Value = 235
$ws.Range("D17").Value = 473
$ws.Range("E17").Value = 24413
$ws.Range("F17").Value = 1.79061063923433
$ws.Range("G17").Value = 41.134326251147
$ws.Range("H17").Value = 4.15669396807273

$ws.Range("B18").Value = "0,6"
$ws.Range("C18").Value = 324
$ws.Range("D18").Value = 582
$ws.Range("E18").Value = 12517
$ws.Range("F18").Value = 2.50888497022094
$ws.Range("G18").Value = 4.996
$ws.Range("H18").Value = 0.63369561369349
